$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.315.45"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "1.708.08"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5312"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2659"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07629"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.571"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").Value = "1.698.85"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "1.943.11"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.65%  "
$ws.Range("D16").Value = "0.0₅8171"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "27.306.58"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("E22").Value = "  -4.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.961"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.43%  "
$ws.Range("E25").Value = "  +7.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1214"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.260"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05391"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.498"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.429"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.643"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.870"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9483"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.872"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.045.99"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8416"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "1.849.29"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4507"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05243"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.38%  "
